$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: advance the date in A1 by one day
$ws.Range("A1").Value = 45309

# Step 2: update prices in column D for rows 33-38
$ws.Range("D33").Value = 186.306
$ws.Range("D34").Value = 148.649
$ws.Range("D35").Value = 138.739
$ws.Range("D36").Value = 319.099
$ws.Range("D37").Value = 227.928
$ws.Range("D38").Value = 200.18
